$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing figures (row 32 / row 33) ---
# KPI "minimize" row: realisasi (C32) drops from 3 to 1 -> deviation worsens
$ws.Range("C32").Value = 1
# KPI "maximize" row: realisasi (D33) drops from 10 to 0 -> deviation worsens
$ws.Range("D33").Value = 0

# --- New "Devisiasi" (deviation) report block, rows 40-42 ---
$ws.Range("C40").Value = "T"
$ws.Range("D40").Value = "C"
$ws.Range("E40").Value = "Score"
$ws.Range("F40").Value = "Devisiasi"
$ws.Range("C40:F40").Font.Bold = $true

$ws.Range("B41").Value = "minimize"
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = 4
$ws.Range("E41").Formula = "=(C41/D41)*100"
$ws.Range("F41").Formula = "=E41-100"
$ws.Range("G41").Value = "%"

$ws.Range("B42").Value = "maximize"
$ws.Range("C42").Value = 10
$ws.Range("D42").Value = 8
$ws.Range("E42").Formula = "=(D42/C42)*100"
$ws.Range("F42").Formula = "=E42-100"
$ws.Range("G42").Value = "%"

# --- View state: scroll / select near the new report block ---
$ws.Activate()
$ws.Range("K39").Select()
